$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '67.464.49'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -1.35%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.765.62'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.44%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '594.46'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -0.14%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '166.09'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.75%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.763.36'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.43%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.518'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -0.30%  '
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -0.12%  '
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +0.31%  '
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -1.39%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '36.21'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +0.54%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.396.68'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +0.44%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.795.53'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +1.30%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '18.41'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +2.75%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '67.453.41'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -1.37%  '
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +0.15%  '
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -0.28%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.03'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -6.78%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '456.49'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -1.93%  '
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +0.01%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.0000154'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +6.56%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '83.19'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -1.57%  '
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -0.63%  '
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -2.17%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.11'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +1.00%  '
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +0.04%  '
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -0.05%  '
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -0.49%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '29.69'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -0.51%  '
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +0.36%  '
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +0.00%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.999'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -0.22%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.717.08'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +0.37%  '
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -0.84%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.33'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -1.19%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.138'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -0.83%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.998'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -0.06%  '
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -1.01%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.00'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +0.01%  '
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +0.01%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '45.34'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +3.17%  '
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -1.72%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '47.04'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +2.38%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '8.33'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -2.72%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '148.41'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +1.08%  '
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -4.43%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '389.66'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +0.20%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '25.58'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -0.40%  '
